# Auto-generated edit script: reassigns species-observation data across rows 6-14
# to match the target workbook state (rows' content permuted + a few blank/comment
# placeholder cells added or removed), per the supplied OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 ---
$ws.Range("A6").Value2 = 111742138
$ws.Range("B6").Value2 = 92683
$ws.Range("D6").Value2 = 'LC'
$ws.Range("E6").Value2 = 2362
$ws.Range("F6").Value2 = 'Blek stjärnmossa'
$ws.Range("G6").Value2 = 'Mnium stellare'
$ws.Range("H6").Value2 = 'Hedw.'
$ws.Range("Q6").Value2 = 331789.362964866
$ws.Range("R6").Value2 = 6626790.418441398

# --- Row 7 ---
$ws.Range("A7").Value2 = 111742077
$ws.Range("B7").Value2 = 78605
$ws.Range("D7").Value2 = 'LC'
$ws.Range("E7").Value2 = 6462
$ws.Range("F7").Value2 = 'Stuplav'
$ws.Range("G7").Value2 = 'Nephroma bellum'
$ws.Range("H7").Value2 = '(Spreng.) Tuck.'

# --- Row 8 ---
$ws.Range("A8").Value2 = 111742184
$ws.Range("B8").Value2 = 93159
$ws.Range("E8").Value2 = 2666
$ws.Range("F8").Value2 = 'Grov fjädermossa'
$ws.Range("G8").Value2 = 'Neckera crispa'
$ws.Range("H8").Value2 = 'Hedw.'
$ws.Range("Q8").Value2 = 331833.6062344447
$ws.Range("R8").Value2 = 6626784.887086328
$ws.Range("AC8").Value2 = 'I bergsbrant'
$ws.Range("L8").NumberFormat = "General"

# --- Row 9 ---
$ws.Range("A9").Value2 = 111742151
$ws.Range("B9").Value2 = 95524
$ws.Range("D9").Value2 = 'LC'
$ws.Range("E9").Value2 = 221944
$ws.Range("F9").Value2 = 'Lopplummer'
$ws.Range("G9").Value2 = 'Huperzia selago'
$ws.Range("H9").Value2 = '(L.) Bernh. ex Schrank & Mart.'
$ws.Range("Q9").Value2 = 331814.6184995985
$ws.Range("R9").Value2 = 6626778.67820756
$ws.Range("AC9").ClearContents() | Out-Null

# --- Row 10 ---
$ws.Range("A10").Value2 = 111742181
$ws.Range("B10").Value2 = 93158
$ws.Range("E10").Value2 = 2667
$ws.Range("F10").Value2 = 'Platt fjädermossa'
$ws.Range("G10").Value2 = 'Neckera complanata'
$ws.Range("H10").Value2 = '(Hedw.) Huebener'
$ws.Range("Q10").Value2 = 331833.6062344447
$ws.Range("R10").Value2 = 6626784.887086328
$ws.Range("AC10").Value2 = 'I bergsbrant'

# --- Row 11 ---
$ws.Range("A11").Value2 = 111742070
$ws.Range("B11").Value2 = 78578
$ws.Range("D11").Value2 = 'NT'
$ws.Range("E11").Value2 = 6458
$ws.Range("F11").Value2 = 'Lunglav'
$ws.Range("G11").Value2 = 'Lobaria pulmonaria'
$ws.Range("H11").Value2 = '(L.) Hoffm.'
$ws.Range("Q11").Value2 = 331735.1116598135
$ws.Range("R11").Value2 = 6626820.629936518
$ws.Range("AC11").Value2 = 'På rönn'
$ws.Range("L11").ClearContents() | Out-Null

# --- Row 12 ---
$ws.Range("A12").Value2 = 111742096
$ws.Range("B12").Value2 = 94125
$ws.Range("D12").Value2 = 'NT'
$ws.Range("E12").Value2 = 54
$ws.Range("F12").Value2 = 'Skogstrappmossa'
$ws.Range("G12").Value2 = 'Anastrophyllum michauxii'
$ws.Range("H12").Value2 = '(F.Weber.) H.Buch'
$ws.Range("Q12").Value2 = 331779.6127968954
$ws.Range("R12").Value2 = 6626798.429951042
$ws.Range("AC12").Value2 = 'På både ved och på lodyta'

# --- Row 13 ---
$ws.Range("A13").Value2 = 111742170
$ws.Range("B13").Value2 = 89369
$ws.Range("E13").Value2 = 5447
$ws.Range("F13").Value2 = 'Vedticka'
$ws.Range("G13").Value2 = 'Fuscoporia viticola'
$ws.Range("H13").Value2 = '(Schwein.) Murrill'
$ws.Range("Q13").Value2 = 331846.7251686137
$ws.Range("R13").Value2 = 6626784.294692003
$ws.Range("L13").ClearContents() | Out-Null
$ws.Range("AC13").ClearContents() | Out-Null

# --- Row 14 ---
$ws.Range("A14").Value2 = 111742101
$ws.Range("B14").Value2 = 94134
$ws.Range("D14").Value2 = 'NT'
$ws.Range("E14").Value2 = 53
$ws.Range("F14").Value2 = 'Vedtrappmossa'
$ws.Range("G14").Value2 = 'Crossocalyx hellerianus'
$ws.Range("H14").Value2 = '(Nees ex Lindenb.) Meyl.'
$ws.Range("Q14").Value2 = 331779.6127968954
$ws.Range("R14").Value2 = 6626798.429951042
$ws.Range("L14").NumberFormat = "General"

Write-Output "Edit complete."
